$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2025-07-21 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-22 Tuesday", 2) | Out-Null

# Update each table cell value by position (row, col) to avoid ambiguity with duplicate values
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "9+84="  # was 63-6=
$t.Cell(1,2).Range.Text = "79+1="  # was 16+36=
$t.Cell(1,3).Range.Text = "56+37="  # was 92-59=
$t.Cell(1,4).Range.Text = "67-38="  # was 33+9=
$t.Cell(1,5).Range.Text = "31+41="  # was 76+19=
$t.Cell(2,1).Range.Text = "78+9="  # was 29-9=
$t.Cell(2,2).Range.Text = "68-29="  # was 86+2=
$t.Cell(2,3).Range.Text = "86+1="  # was 94-31=
$t.Cell(2,4).Range.Text = "10+76="  # was 13-5=
$t.Cell(2,5).Range.Text = "58-34="  # was 2+60=
$t.Cell(3,1).Range.Text = "97-80="  # was 22+44=
$t.Cell(3,2).Range.Text = "83-5="  # was 48-15=
$t.Cell(3,3).Range.Text = "68-46="  # was 78-47=
$t.Cell(3,4).Range.Text = "72-71="  # was 16+37=
$t.Cell(3,5).Range.Text = "99-55="  # was 92-22=
$t.Cell(4,1).Range.Text = "63-42="  # was 96-60=
$t.Cell(4,2).Range.Text = "20+67="  # was 50-49=
$t.Cell(4,3).Range.Text = "48+40="  # was 46+2=
$t.Cell(4,4).Range.Text = "72+25="  # was 12+45=
$t.Cell(4,5).Range.Text = "79+11="  # was 79-35=
$t.Cell(5,1).Range.Text = "76-16="  # was 37+49=
$t.Cell(5,2).Range.Text = "16-5="  # was 92-40=
$t.Cell(5,3).Range.Text = "4+30="  # was 25+46=
$t.Cell(5,4).Range.Text = "12-8="  # was 30+52=
$t.Cell(5,5).Range.Text = "28+51="  # was 76+19=
$t.Cell(6,1).Range.Text = "42-27="  # was 44-21=
$t.Cell(6,2).Range.Text = "87-2="  # was 40+3=
$t.Cell(6,3).Range.Text = "74-57="  # was 28+9=
$t.Cell(6,4).Range.Text = "19+51="  # was 56-0=
$t.Cell(6,5).Range.Text = "21+6="  # was 26+42=
$t.Cell(7,1).Range.Text = "39+29="  # was 82+5=
$t.Cell(7,2).Range.Text = "99-15="  # was 56-24=
$t.Cell(7,3).Range.Text = "2+92="  # was 82-30=
$t.Cell(7,4).Range.Text = "3+62="  # was 13+11=
$t.Cell(7,5).Range.Text = "82+17="  # was 45-16=
$t.Cell(8,1).Range.Text = "27+69="  # was 64+26=
$t.Cell(8,2).Range.Text = "0+94="  # was 5+65=
$t.Cell(8,3).Range.Text = "66+31="  # was 16+19=
$t.Cell(8,4).Range.Text = "87+11="  # was 34-7=
$t.Cell(8,5).Range.Text = "42+33="  # was 27+7=
$t.Cell(9,1).Range.Text = "89-13="  # was 88-50=
$t.Cell(9,2).Range.Text = "2+95="  # was 14+79=
$t.Cell(9,3).Range.Text = "37+12="  # was 64+35=
$t.Cell(9,4).Range.Text = "70-3="  # was 87-67=
$t.Cell(9,5).Range.Text = "29+51="  # was 43-40=
$t.Cell(10,1).Range.Text = "52-15="  # was 78-31=
$t.Cell(10,2).Range.Text = "67-46="  # was 89-57=
$t.Cell(10,3).Range.Text = "6+38="  # was 64-10=
$t.Cell(10,4).Range.Text = "2+2="  # was 53+36=
$t.Cell(10,5).Range.Text = "31+40="  # was 13+9=
$t.Cell(11,1).Range.Text = "40+27="  # was 93-38=
$t.Cell(11,2).Range.Text = "48+2="  # was 20+25=
$t.Cell(11,3).Range.Text = "15+76="  # was 8+33=
$t.Cell(11,4).Range.Text = "82-5="  # was 43-30=
$t.Cell(11,5).Range.Text = "56-33="  # was 59-1=
$t.Cell(12,1).Range.Text = "13+7="  # was 53+41=
$t.Cell(12,2).Range.Text = "72-4="  # was 47+32=
$t.Cell(12,3).Range.Text = "40+52="  # was 68-31=
$t.Cell(12,4).Range.Text = "79-0="  # was 23+23=
$t.Cell(12,5).Range.Text = "33-6="  # was 76+9=
$t.Cell(13,1).Range.Text = "72-43="  # was 29+62=
$t.Cell(13,2).Range.Text = "18+49="  # was 73-38=
$t.Cell(13,3).Range.Text = "50+40="  # was 53+31=
$t.Cell(13,4).Range.Text = "21-6="  # was 39+19=
$t.Cell(13,5).Range.Text = "96-57="  # was 49+36=
$t.Cell(14,1).Range.Text = "60-5="  # was 24+31=
$t.Cell(14,2).Range.Text = "8+20="  # was 64+34=
$t.Cell(14,3).Range.Text = "77-53="  # was 91-60=
$t.Cell(14,4).Range.Text = "65+22="  # was 23+61=
$t.Cell(14,5).Range.Text = "60-23="  # was 18+15=
$t.Cell(15,1).Range.Text = "83-79="  # was 70-12=
$t.Cell(15,2).Range.Text = "70-1="  # was 31+47=
$t.Cell(15,3).Range.Text = "13+51="  # was 12+52=
$t.Cell(15,4).Range.Text = "49-21="  # was 27+32=
$t.Cell(15,5).Range.Text = "40+56="  # was 72+5=
$t.Cell(16,1).Range.Text = "87-86="  # was 33+25=
$t.Cell(16,2).Range.Text = "2+1="  # was 66-53=
$t.Cell(16,3).Range.Text = "41-40="  # was 52+6=
$t.Cell(16,4).Range.Text = "89-30="  # was 18-12=
$t.Cell(16,5).Range.Text = "72-62="  # was 76-14=
$t.Cell(17,1).Range.Text = "98-51="  # was 51+19=
$t.Cell(17,2).Range.Text = "41-23="  # was 50-28=
$t.Cell(17,3).Range.Text = "54+38="  # was 80-29=
$t.Cell(17,4).Range.Text = "6+43="  # was 83-74=
$t.Cell(17,5).Range.Text = "48-29="  # was 6+42=
$t.Cell(18,1).Range.Text = "80-73="  # was 32+64=
$t.Cell(18,2).Range.Text = "4+72="  # was 4+78=
$t.Cell(18,3).Range.Text = "93-4="  # was 18+36=
$t.Cell(18,4).Range.Text = "83-17="  # was 64-51=
$t.Cell(18,5).Range.Text = "71+26="  # was 95-71=
$t.Cell(19,1).Range.Text = "97-97="  # was 38-14=
$t.Cell(19,2).Range.Text = "55+3="  # was 68+4=
$t.Cell(19,3).Range.Text = "58+1="  # was 7+40=
$t.Cell(19,4).Range.Text = "99-43="  # was 60-18=
$t.Cell(19,5).Range.Text = "91+4="  # was 41+51=
$t.Cell(20,1).Range.Text = "30+23="  # was 49-26=
$t.Cell(20,2).Range.Text = "46+16="  # was 35+32=
$t.Cell(20,3).Range.Text = "92-90="  # was 45+35=
$t.Cell(20,4).Range.Text = "41+14="  # was 10-0=
$t.Cell(20,5).Range.Text = "94-14="  # was 32-24=
